# Update "想去人数" (people interested) counts in column F
# for both the "展览" sheet and the "全部类型" sheet (which mirrors it).

$wb = $excel.ActiveWorkbook

# Row -> new F value
$updates = @{
    7  = 1214
    8  = 1501
    10 = 373
    12 = 136
    17 = 291
    18 = 316
    19 = 1704
    25 = 332
    26 = 4101
    30 = 1069
    33 = 443
    35 = 206
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
